$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.934.70"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "2.811.74"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.05"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.38"
$ws.Range("E6").Value = "  +5.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  +9.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.42"
$ws.Range("E10").Value = "  +4.56%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0840"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.00"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.80"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("D15").Value = "3.247.92"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "2.810.37"
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.948"
$ws.Range("E17").Value = "  +4.13%  "
$ws.Range("D18").Value = "51.871.20"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.65"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  +8.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.63"
$ws.Range("E21").Value = "  +6.62%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  +3.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.38"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.90"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.22"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  +14.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.40"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.57"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.14"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0459"
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0895"
$ws.Range("E35").Value = "  +9.26%  "
$ws.Range("E36").Value = "  +8.27%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.96"
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("E39").Value = "  +5.23%  "
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("E41").Value = "  +3.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.54"
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.21"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.03"
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.50"
$ws.Range("E46").Value = "  +10.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("E47").Value = "  +6.74%  "
$ws.Range("D48").Value = "2.108.14"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.966"
$ws.Range("E49").Value = "  +7.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.50"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("E51").Value = "  +8.25%  "
